$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new value as a literal-text formula first (so Excel does not
# "smart" coerce the comma-separated digit string into a huge number), then
# flatten the formula down to a plain value in place via Copy/PasteSpecial
# (values only) so the stored cell ends up as literal text - matching how
# the source file stores it (t="s", same original style untouched).
$ws.Range("C2").Formula = '="111757,222222,333333,444444"'
$ws.Range("C2").Copy()
$ws.Range("C2").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Range("C2").Select()
